$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Bump the cached "today" date field (datetimeFigureOut) from 4/1/2024 to
#    4/3/2024 everywhere it is rendered from: the Slide Master, every Custom
#    Layout, and the Notes Master.
# ---------------------------------------------------------------------------
function Update-DateField {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -and $shp.Name -like "*Date*") {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "4/1/2024") {
                $tr.Text = "4/3/2024"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateField $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateField $layout.Shapes
}

Update-DateField $p.NotesMaster.Shapes

# ---------------------------------------------------------------------------
# 2) Slide 1 subtitle: "NCSU scRNA Workshop, 2024" / "Bruce Corliss, PhD and
#    Allison Dickey, PhD" -> single line "NCSU scRNA-Seq Workshop, 2024"
#    (drop the byline paragraph entirely).
# ---------------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
for ($i = 1; $i -le $slide1.Shapes.Count; $i++) {
    $shp = $slide1.Shapes.Item($i)
    if ($shp.HasTextFrame -and $shp.Name -eq "Subtitle 2") {
        $tr = $shp.TextFrame.TextRange

        # Turn the run " Workshop, 2024" into "-Seq Workshop, 2024", leaving
        # the preceding "NCSU " / "scRNA" runs untouched.
        $prefix = "NCSU scRNA"
        $oldTail = " Workshop, 2024"
        $newTail = "-Seq Workshop, 2024"
        $tailRange = $tr.Characters($prefix.Length + 1, $oldTail.Length)
        if ($tailRange.Text -eq $oldTail) {
            $tailRange.Text = $newTail
        }

        # Remove the second paragraph ("Bruce Corliss, PhD and Allison
        # Dickey, PhD") in its entirety.
        $newFirstParaLen = ($prefix + $newTail).Length
        $restStart = $newFirstParaLen + 2
        if ($restStart -le $tr.Length) {
            $restRange = $tr.Characters($restStart, $tr.Length)
            $restRange.Delete()
        }
    }
}
